$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.734.91'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '2.022.66'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.74'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.54'
$ws.Range("E7").Value = '  +3.37%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0812'
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.52'
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.321.86'
$ws.Range("E13").Value = '  -1.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.92'
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").Value = '2.041.82'
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '37.774.05'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.00'
$ws.Range("E19").Value = '  -4.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.53'
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = '0.0₃0818'
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.32'
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("E24").Value = '  -2.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.67'
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.17'
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.130'
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.83'
$ws.Range("E29").Value = '  -1.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.28'
$ws.Range("E30").Value = '  -6.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.42'
$ws.Range("E32").Value = '  -3.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.05'
$ws.Range("E33").Value = '  +3.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0600'
$ws.Range("E34").Value = '  -2.74%  '
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.33'
$ws.Range("E36").Value = '  +5.74%  '
$ws.Range("E37").Value = '  -5.91%  '
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("D40").Value = '1.534.95'
$ws.Range("E40").Value = '  +3.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0216'
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.38'
$ws.Range("E42").Value = '  -2.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.81'
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.47'
$ws.Range("E44").Value = '  -1.83%  '
$ws.Range("E45").Value = '  -3.84%  '
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.95'
$ws.Range("E47").Value = '  -3.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.95'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("E49").Value = '  -2.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.09'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").Value = '2.213.92'
$ws.Range("E51").Value = '  -1.49%  '
